# Refresh the crypto price/volume table (columns D "Price" and E "Volume(1h)")
# with the latest scraped values, mirroring the GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.200.88'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '1.601.64'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'0.9999"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = "'302.90"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('D7').Value = "'0.3779"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = "'51.66"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.16%  '
$ws.Range('D9').Value = "'0.3611"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.12%  '
$ws.Range('D10').Value = "'1.264"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').Value = "'1.001"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = "'0.08120"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').Value = "'22.57"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').Value = "'6.582"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').Value = "'7.385"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('E16').Value = '  -0.62%  '
$ws.Range('D17').Value = '1.601.82'
$ws.Range('D18').Value = "'93.60"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.21%  '
$ws.Range('D19').Value = "'0.06862"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').Value = "'18.02"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').Value = "'0.9995"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').Value = "'12.94"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').Value = '23.193.06'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('E25').Value = '  +2.34%  '
$ws.Range('D26').Value = "'2.995"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.82%  '
$ws.Range('D27').Value = "'21.17"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('D28').Value = "'150.09"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').Value = "'5.235"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.89%  '
$ws.Range('D30').Value = "'133.60"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('D31').Value = "'2.417"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('D32').Value = "'6.816"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('D33').Value = '1.779.13'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').Value = "'0.9785"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.05%  '
$ws.Range('D35').Value = "'0.07544"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.97%  '
$ws.Range('D36').Value = "'10.33"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.82%  '
$ws.Range('D37').Value = "'0.02723"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('D38').Value = "'6.141"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('D39').Value = "'0.2501"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.72%  '
$ws.Range('D40').Value = "'0.08795"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('D41').Value = "'0.7096"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').Value = "'1.359"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.94%  '
$ws.Range('D43').Value = "'12.42"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.29%  '
$ws.Range('D44').Value = "'15.40"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('D45').Value = "'0.6542"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.28%  '
$ws.Range('D46').Value = "'2.307"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').Value = "'4.015"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.02%  '
$ws.Range('D48').Value = "'132.31"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').Value = "'0.07961"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('D50').Value = "'1.203"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.42%  '
$ws.Range('D51').Value = "'1.231"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.92%  '
